$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source survey data was re-parsed; rows were re-sorted by donor_id and
# several rows worth of newly-parsed donors were added, growing the table from
# 8 data rows (A1:E9) to 23 data rows (A1:E24).

# Insert 15 blank rows (new formatting is inherited from the row above, exactly
# like typing rows in the Excel UI) so the sheet grows from 9 to 24 total rows.
$ws.Range("A9:A23").EntireRow.Insert()

# Re-populate every data row (2-24) with the freshly re-parsed values.
$ws.Cells.Item(2,1).Value = "0ce5dd49"
$ws.Cells.Item(2,2).Value = 2.339066339066339
$ws.Cells.Item(2,3).Formula = "=""2-3"""
$ws.Cells.Item(2,3).Copy()
$ws.Cells.Item(2,3).PasteSpecial(-4163)
$ws.Cells.Item(2,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(2,5).Value = 45854.63724849415

$ws.Cells.Item(3,1).Value = "2c1001cb"
$ws.Cells.Item(3,2).Value = 2.848648648648649
$ws.Cells.Item(3,3).Formula = "=""2-3"""
$ws.Cells.Item(3,3).Copy()
$ws.Cells.Item(3,3).PasteSpecial(-4163)
$ws.Cells.Item(3,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(3,5).Value = 45854.63724849415

$ws.Cells.Item(4,1).Value = "37cc37bf"
$ws.Cells.Item(4,2).Value = 1.295774647887324
$ws.Cells.Item(4,3).Formula = "=""1"""
$ws.Cells.Item(4,3).Copy()
$ws.Cells.Item(4,3).PasteSpecial(-4163)
$ws.Cells.Item(4,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(4,5).Value = 45854.63724849415

$ws.Cells.Item(5,1).Value = "43faa0b9"
$ws.Cells.Item(5,2).Value = 1.738983050847458
$ws.Cells.Item(5,3).Formula = "=""1"""
$ws.Cells.Item(5,3).Copy()
$ws.Cells.Item(5,3).PasteSpecial(-4163)
$ws.Cells.Item(5,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(5,5).Value = 45854.63724849415

$ws.Cells.Item(6,1).Value = "4abe3e88"
$ws.Cells.Item(6,2).Value = 1.321167883211679
$ws.Cells.Item(6,3).Formula = "=""1"""
$ws.Cells.Item(6,3).Copy()
$ws.Cells.Item(6,3).PasteSpecial(-4163)
$ws.Cells.Item(6,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(6,5).Value = 45854.63724849415

$ws.Cells.Item(7,1).Value = "50164f59"
$ws.Cells.Item(7,2).Value = 1.963768115942029
$ws.Cells.Item(7,3).Formula = "=""1"""
$ws.Cells.Item(7,3).Copy()
$ws.Cells.Item(7,3).PasteSpecial(-4163)
$ws.Cells.Item(7,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(7,5).Value = 45854.63724849415

$ws.Cells.Item(8,1).Value = "5cf70f79"
$ws.Cells.Item(8,2).Value = 1.541666666666667
$ws.Cells.Item(8,3).Formula = "=""1"""
$ws.Cells.Item(8,3).Copy()
$ws.Cells.Item(8,3).PasteSpecial(-4163)
$ws.Cells.Item(8,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(8,5).Value = 45854.63724849415

$ws.Cells.Item(9,1).Value = "5da96769"
$ws.Cells.Item(9,2).Value = 1.5
$ws.Cells.Item(9,3).Formula = "=""1"""
$ws.Cells.Item(9,3).Copy()
$ws.Cells.Item(9,3).PasteSpecial(-4163)
$ws.Cells.Item(9,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(9,5).Value = 45854.63724849415

$ws.Cells.Item(10,1).Value = "6ca3e2f6"
$ws.Cells.Item(10,2).Value = 1.1
$ws.Cells.Item(10,3).Formula = "=""1"""
$ws.Cells.Item(10,3).Copy()
$ws.Cells.Item(10,3).PasteSpecial(-4163)
$ws.Cells.Item(10,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(10,5).Value = 45854.63724849415

$ws.Cells.Item(11,1).Value = "790a4fcb"
$ws.Cells.Item(11,2).Value = 1.3
$ws.Cells.Item(11,3).Formula = "=""1"""
$ws.Cells.Item(11,3).Copy()
$ws.Cells.Item(11,3).PasteSpecial(-4163)
$ws.Cells.Item(11,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(11,5).Value = 45854.63724849415

$ws.Cells.Item(12,1).Value = "802cc63a"
$ws.Cells.Item(12,2).Value = 1.525423728813559
$ws.Cells.Item(12,3).Formula = "=""1"""
$ws.Cells.Item(12,3).Copy()
$ws.Cells.Item(12,3).PasteSpecial(-4163)
$ws.Cells.Item(12,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(12,5).Value = 45854.63724849415

$ws.Cells.Item(13,1).Value = "85c3ea4d"
$ws.Cells.Item(13,2).Value = 3.57396449704142
$ws.Cells.Item(13,3).Formula = "=""2-3"""
$ws.Cells.Item(13,3).Copy()
$ws.Cells.Item(13,3).PasteSpecial(-4163)
$ws.Cells.Item(13,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(13,5).Value = 45854.63724849415

$ws.Cells.Item(14,1).Value = "942dfafb"
$ws.Cells.Item(14,2).Value = 2.426373626373626
$ws.Cells.Item(14,3).Formula = "=""2-3"""
$ws.Cells.Item(14,3).Copy()
$ws.Cells.Item(14,3).PasteSpecial(-4163)
$ws.Cells.Item(14,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(14,5).Value = 45854.63724849415

$ws.Cells.Item(15,1).Value = "9bc6ba8c"
$ws.Cells.Item(15,2).Value = 1
$ws.Cells.Item(15,3).Formula = "=""1"""
$ws.Cells.Item(15,3).Copy()
$ws.Cells.Item(15,3).PasteSpecial(-4163)
$ws.Cells.Item(15,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(15,5).Value = 45854.63724849415

$ws.Cells.Item(16,1).Value = "a2d65af2"
$ws.Cells.Item(16,2).Value = 1.1
$ws.Cells.Item(16,3).Formula = "=""1"""
$ws.Cells.Item(16,3).Copy()
$ws.Cells.Item(16,3).PasteSpecial(-4163)
$ws.Cells.Item(16,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(16,5).Value = 45854.63724849415

$ws.Cells.Item(17,1).Value = "a46f1771"
$ws.Cells.Item(17,2).Value = 1.35
$ws.Cells.Item(17,3).Formula = "=""1"""
$ws.Cells.Item(17,3).Copy()
$ws.Cells.Item(17,3).PasteSpecial(-4163)
$ws.Cells.Item(17,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(17,5).Value = 45854.63724849415

$ws.Cells.Item(18,1).Value = "ad58f9da"
$ws.Cells.Item(18,2).Value = 1.161016949152542
$ws.Cells.Item(18,3).Formula = "=""1"""
$ws.Cells.Item(18,3).Copy()
$ws.Cells.Item(18,3).PasteSpecial(-4163)
$ws.Cells.Item(18,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(18,5).Value = 45854.63724849415

$ws.Cells.Item(19,1).Value = "c7d9a301"
$ws.Cells.Item(19,2).Value = 1.578313253012048
$ws.Cells.Item(19,3).Formula = "=""1"""
$ws.Cells.Item(19,3).Copy()
$ws.Cells.Item(19,3).PasteSpecial(-4163)
$ws.Cells.Item(19,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(19,5).Value = 45854.63724849415

$ws.Cells.Item(20,1).Value = "ce8732ff"
$ws.Cells.Item(20,2).Value = 1.529801324503311
$ws.Cells.Item(20,3).Formula = "=""1"""
$ws.Cells.Item(20,3).Copy()
$ws.Cells.Item(20,3).PasteSpecial(-4163)
$ws.Cells.Item(20,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(20,5).Value = 45854.63724849415

$ws.Cells.Item(21,1).Value = "d6f1d567"
$ws.Cells.Item(21,2).Value = 1.461538461538461
$ws.Cells.Item(21,3).Formula = "=""1"""
$ws.Cells.Item(21,3).Copy()
$ws.Cells.Item(21,3).PasteSpecial(-4163)
$ws.Cells.Item(21,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(21,5).Value = 45854.63724849415

$ws.Cells.Item(22,1).Value = "da9326c9"
$ws.Cells.Item(22,2).Value = 1.857142857142857
$ws.Cells.Item(22,3).Formula = "=""1"""
$ws.Cells.Item(22,3).Copy()
$ws.Cells.Item(22,3).PasteSpecial(-4163)
$ws.Cells.Item(22,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(22,5).Value = 45854.63724849415

$ws.Cells.Item(23,1).Value = "e09ca7bf"
$ws.Cells.Item(23,2).Value = 2.426373626373626
$ws.Cells.Item(23,3).Formula = "=""2-3"""
$ws.Cells.Item(23,3).Copy()
$ws.Cells.Item(23,3).PasteSpecial(-4163)
$ws.Cells.Item(23,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(23,5).Value = 45854.63724849415

$ws.Cells.Item(24,1).Value = "ef53a641"
$ws.Cells.Item(24,2).Value = 2.339066339066339
$ws.Cells.Item(24,3).Formula = "=""2-3"""
$ws.Cells.Item(24,3).Copy()
$ws.Cells.Item(24,3).PasteSpecial(-4163)
$ws.Cells.Item(24,4).Value = "q02_sessions_per_active_day"
$ws.Cells.Item(24,5).Value = 45854.63724849415

$excel.CutCopyMode = $false
